$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(17, 18),
    @(41, 42),
    @(58, 59),
    @(69, 70),
    @(73, 74),
    @(78, 79),
    @(91, 92),
    @(103, 104),
    @(135, 136),
    @(173, 174),
    @(190, 191),
    @(223, 224)
)

$firstCol = 2   # column B
$lastCol = 30   # column AD

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
